$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by
# Excel's smart-entry parsing (e.g. "1.00" becomes 1, "288.33" becomes a
# float that round-trips with binary noise, "0.0497" becomes scientific
# notation) are forced to Text format first so the literal string survives,
# matching the source data's inline-string cells.

$ws.Range("D2").Value = "57.815.03"
$ws.Range("E2").Value = "  -4.18%  "
$ws.Range("D3").Value = "2.289.29"
$ws.Range("E3").Value = "  -5.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.15"
$ws.Range("E5").Value = "  -4.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.52"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -3.36%  "
$ws.Range("D9").Value = "2.288.23"
$ws.Range("E9").Value = "  -5.05%  "
$ws.Range("E10").Value = "  -5.85%  "
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.53"
$ws.Range("E14").Value = "  -4.86%  "
$ws.Range("D15").Value = "2.697.84"
$ws.Range("E15").Value = "  -5.01%  "
$ws.Range("D16").Value = "57.821.08"
$ws.Range("E16").Value = "  -4.03%  "
$ws.Range("E17").Value = "  -5.10%  "
$ws.Range("D18").Value = "2.274.88"
$ws.Range("E18").Value = "  -5.51%  "
$ws.Range("E19").Value = "  -5.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.24"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.70"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("E22").Value = "  -6.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.90"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("E25").Value = "  -4.98%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -6.58%  "
$ws.Range("E28").Value = "  -5.26%  "
$ws.Range("E29").Value = "  -5.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.32"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").Value = "0.0₃0722"
$ws.Range("E31").Value = "  -6.55%  "
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("E33").Value = "  -6.68%  "
$ws.Range("E34").Value = "  -6.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.66"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -7.91%  "
$ws.Range("E39").Value = "  -7.03%  "
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("E41").Value = "  -7.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.98"
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "288.33"
$ws.Range("E43").Value = "  -11.05%  "
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("E45").Value = "  -2.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0497"
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.23"
$ws.Range("E48").Value = "  -8.54%  "
$ws.Range("E49").Value = "  -5.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.95"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  -0.44%  "
